$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are plain numeric-looking text in the source data
# (e.g. "30.107.75", "1.001") - not real numbers. A leading apostrophe forces
# Excel to store them as text instead of coercing to a number; resetting the
# style back to Normal afterwards clears the "quote prefix" flag Excel sets,
# matching the original cell formatting (no explicit style).

$ws.Range("D2").Value = "'30.137.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.69%  "

$ws.Range("D3").Value = "'1.918.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.56%  "

$ws.Range("D5").Value = "'329.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.32%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.80%  "

$ws.Range("D7").Value = "'0.5211"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.39%  "

$ws.Range("D8").Value = "'0.4082"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.60%  "

$ws.Range("D9").Value = "'0.08523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.05%  "

$ws.Range("D10").Value = "'42.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.02%  "

$ws.Range("D11").Value = "'1.123"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "'22.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.33%  "

$ws.Range("D13").Value = "'6.420"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.52%  "

$ws.Range("D14").Value = "'1.918.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.72%  "

$ws.Range("D15").Value = "'7.388"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.84%  "

$ws.Range("E16").Value = "  -0.87%  "

$ws.Range("D17").Value = "'95.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.17%  "

$ws.Range("D18").Value = "'0.00001114"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.40%  "

$ws.Range("D19").Value = "'0.06699"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("D20").Value = "'18.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.26%  "

$ws.Range("E21").Value = "  -0.71%  "

$ws.Range("D22").Value = "'6.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.63%  "

$ws.Range("D23").Value = "'30.141.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.50%  "

$ws.Range("D24").Value = "'11.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.43%  "

$ws.Range("D25").Value = "'2.216"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("D26").Value = "'2.139.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.64%  "

$ws.Range("D27").Value = "'160.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.00%  "

$ws.Range("D28").Value = "'21.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.69%  "

$ws.Range("D29").Value = "'2.433"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.59%  "

$ws.Range("D30").Value = "'128.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.66%  "

$ws.Range("D31").Value = "'1.079"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.57%  "

$ws.Range("E32").Value = "  +2.53%  "

$ws.Range("D33").Value = "'6.038"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.38%  "

$ws.Range("D34").Value = "'3.645"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.78%  "

$ws.Range("D35").Value = "'0.02492"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.56%  "

$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("D37").Value = "'0.2204"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.12%  "

$ws.Range("D38").Value = "'1.229"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.22%  "

$ws.Range("D39").Value = "'5.183"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.73%  "

$ws.Range("D40").Value = "'8.884"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.32%  "

$ws.Range("D41").Value = "'0.6543"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.97%  "

$ws.Range("D42").Value = "'11.63"
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").Value = "'1.243"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.72%  "

$ws.Range("D44").Value = "'0.6144"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.56%  "

$ws.Range("D45").Value = "'13.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("D46").Value = "'3.755"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.97%  "

$ws.Range("D47").Value = "'2.078"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.73%  "

$ws.Range("D48").Value = "'1.243"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.52%  "

$ws.Range("D49").Value = "'124.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.65%  "

$ws.Range("D50").Value = "'1.163"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.42%  "

$ws.Range("D51").Value = "'79.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.68%  "
